$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2072.8696
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2072.8696
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6218.6088
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6554.6088
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H62").Value = 9809108
$ws.Range("I62").Value = 11116449
$ws.Range("J62").Value = 4049.5
$ws.Range("K62").Value = 11116449
$ws.Range("L62").Value = 4049.5
$ws.Range("M62").Value = -11115825
$ws.Range("N62").Value = -5297.5
$ws.Range("H65").Value = 9809108
$ws.Range("I65").Value = 11116449
$ws.Range("J65").Value = 4049.5
$ws.Range("K65").Value = 55582245
$ws.Range("L65").Value = 20247.5
$ws.Range("M65").Value = -55579125
$ws.Range("N65").Value = -26487.5
$ws.Range("H138").Value = 4279.8057
$ws.Range("I138").Value = 2120.0588
$ws.Range("J138").Value = 6212.2104
$ws.Range("K138").Value = 6360.176399999999
$ws.Range("L138").Value = 18636.6312
$ws.Range("M138").Value = -1220.176399999999
$ws.Range("N138").Value = -28916.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2317315.2
$ws.Range("I32").Value = 2283.152
$ws.Range("K32").Value = 2283.152
$ws.Range("M32").Value = -1996.152
$ws.Range("H61").Value = 4222.7383
$ws.Range("I61").Value = 3420.6
$ws.Range("K61").Value = 3420.6
$ws.Range("M61").Value = -3208.6
$ws.Range("H74").Value = 4737
$ws.Range("I74").Value = 4583.467
$ws.Range("J74").Value = 5197.6
$ws.Range("K74").Value = 4583.467
$ws.Range("L74").Value = 5197.6
$ws.Range("M74").Value = -3709.467
$ws.Range("N74").Value = -6945.6
$ws.Range("H77").Value = 4737
$ws.Range("I77").Value = 4583.467
$ws.Range("J77").Value = 5197.6
$ws.Range("K77").Value = 22917.335
$ws.Range("L77").Value = 25988
$ws.Range("M77").Value = -18549.335
$ws.Range("N77").Value = -34724
$ws.Range("H110").Value = 4923.5757
$ws.Range("I110").Value = 2422.2727
$ws.Range("K110").Value = 2422.2727
$ws.Range("M110").Value = -377.2727
$ws.Range("H122").Value = 3474.2
$ws.Range("I122").Value = 2882.1667
$ws.Range("J122").Value = 4362.25
$ws.Range("K122").Value = 8646.500100000001
$ws.Range("L122").Value = 13086.75
$ws.Range("M122").Value = -6196.500100000001
$ws.Range("N122").Value = -17986.75
$ws.Range("H132").Value = 848667.6
$ws.Range("I132").Value = 960465.0600000001
$ws.Range("K132").Value = 2881395.18
$ws.Range("M132").Value = -2878865.18
$ws.Range("H136").Value = 4222.7383
$ws.Range("I136").Value = 3420.6
$ws.Range("K136").Value = 10261.8
$ws.Range("M136").Value = -7711.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3410.6428
$ws.Range("I105").Value = 3994
$ws.Range("J105").Value = 1952.25
$ws.Range("K105").Value = 3994
$ws.Range("L105").Value = 1952.25
$ws.Range("M105").Value = -2247
$ws.Range("N105").Value = -5446.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 567.5
$ws.Range("I19").Value = 366.33334
$ws.Range("J19").Value = 768.6667
$ws.Range("K19").Value = 366.33334
$ws.Range("L19").Value = 768.6667
$ws.Range("M19").Value = -196.33334
$ws.Range("N19").Value = -1108.6667
$ws.Range("H24").Value = 567.5
$ws.Range("I24").Value = 366.33334
$ws.Range("J24").Value = 768.6667
$ws.Range("K24").Value = 366.33334
$ws.Range("L24").Value = 768.6667
$ws.Range("M24").Value = -196.33334
$ws.Range("N24").Value = -1108.6667
$ws.Range("H31").Value = 4761
$ws.Range("I31").Value = 1295
$ws.Range("K31").Value = 1295
$ws.Range("M31").Value = -1000
$ws.Range("H34").Value = 4761
$ws.Range("I34").Value = 1295
$ws.Range("K34").Value = 1295
$ws.Range("M34").Value = -1093
$ws.Range("H50").Value = 38961
$ws.Range("J50").Value = 38961
$ws.Range("L50").Value = 38961
$ws.Range("N50").Value = -40211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 13748
$ws.Range("J64").Value = 12882.538
$ws.Range("L64").Value = 38647.614
$ws.Range("N64").Value = -39187.614
$ws.Range("H67").Value = 13748
$ws.Range("J67").Value = 12882.538
$ws.Range("L67").Value = 38647.614
$ws.Range("N67").Value = -40519.614
$ws.Range("H131").Value = 48151172
$ws.Range("I131").Value = 35558716
$ws.Range("J131").Value = 111113460
$ws.Range("K131").Value = 106676148
$ws.Range("L131").Value = 333340380
$ws.Range("M131").Value = -106671108
$ws.Range("N131").Value = -333350460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9046.299999999999
$ws.Range("I70").Value = 9308.125
$ws.Range("K70").Value = 9308.125
$ws.Range("M70").Value = -9038.125
$ws.Range("H73").Value = 9046.299999999999
$ws.Range("I73").Value = 9308.125
$ws.Range("K73").Value = 9308.125
$ws.Range("M73").Value = -8372.125
$ws.Range("H80").Value = 3928.0264
$ws.Range("J80").Value = 4442.3335
$ws.Range("L80").Value = 4442.3335
$ws.Range("N80").Value = -6438.3335
$ws.Range("H83").Value = 3928.0264
$ws.Range("J83").Value = 4442.3335
$ws.Range("L83").Value = 22211.6675
$ws.Range("N83").Value = -32195.6675
$ws.Range("H92").Value = 7916.6665
$ws.Range("J92").Value = 7916.6665
$ws.Range("L92").Value = 7916.6665
$ws.Range("N92").Value = -11660.6665
$ws.Range("H122").Value = 4729.857
$ws.Range("I122").Value = 3361.0908
$ws.Range("J122").Value = 6235.5
$ws.Range("K122").Value = 10083.2724
$ws.Range("L122").Value = 18706.5
$ws.Range("M122").Value = -7633.2724
$ws.Range("N122").Value = -23606.5
$ws.Range("H132").Value = 38466444
$ws.Range("I132").Value = 62505252
$ws.Range("J132").Value = 4349.8
$ws.Range("K132").Value = 187515756
$ws.Range("L132").Value = 13049.4
$ws.Range("M132").Value = -187513226
$ws.Range("N132").Value = -18109.4
$ws.Range("H133").Value = 88500
$ws.Range("J133").Value = 88500
$ws.Range("L133").Value = 88500
$ws.Range("N133").Value = -98620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4761.1577
$ws.Range("I55").Value = 4120
$ws.Range("K55").Value = 4120
$ws.Range("M55").Value = -3947
$ws.Range("H68").Value = 3578.2173
$ws.Range("I68").Value = 2841.0588
$ws.Range("J68").Value = 5666.8335
$ws.Range("K68").Value = 2841.0588
$ws.Range("L68").Value = 5666.8335
$ws.Range("M68").Value = -2092.0588
$ws.Range("N68").Value = -7164.8335
$ws.Range("H71").Value = 3578.2173
$ws.Range("I71").Value = 2841.0588
$ws.Range("J71").Value = 5666.8335
$ws.Range("K71").Value = 14205.294
$ws.Range("L71").Value = 28334.1675
$ws.Range("M71").Value = -10461.294
$ws.Range("N71").Value = -35822.1675
$ws.Range("H100").Value = 2664.1155
$ws.Range("I100").Value = 3212.2144
$ws.Range("J100").Value = 2024.6666
$ws.Range("K100").Value = 3212.2144
$ws.Range("L100").Value = 2024.6666
$ws.Range("M100").Value = -2671.2144
$ws.Range("N100").Value = -3106.6666
$ws.Range("H123").Value = 9500
$ws.Range("J123").Value = 9500
$ws.Range("L123").Value = 9500
$ws.Range("N123").Value = -19300
$ws.Range("H136").Value = 45462756
$ws.Range("I136").Value = 100009110
$ws.Range("J136").Value = 7459.5835
$ws.Range("K136").Value = 300027330
$ws.Range("L136").Value = 22378.7505
$ws.Range("M136").Value = -300024780
$ws.Range("N136").Value = -27478.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 15210
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 15210
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 15210
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -15826
$ws.Range("H107").Value = 2422.35
$ws.Range("I107").Value = 1883.8572
$ws.Range("J107").Value = 3678.8333
$ws.Range("K107").Value = 5651.571599999999
$ws.Range("L107").Value = 11036.4999
$ws.Range("M107").Value = -3731.571599999999
$ws.Range("N107").Value = -14876.4999
$ws.Range("H122").Value = 11609.5
$ws.Range("I122").Value = 3978.647
$ws.Range("J122").Value = 30141.572
$ws.Range("K122").Value = 11935.941
$ws.Range("L122").Value = 90424.716
$ws.Range("M122").Value = -9485.940999999999
$ws.Range("N122").Value = -95324.716
$ws.Range("H136").Value = 15162375
$ws.Range("I136").Value = 18527350
$ws.Range("J136").Value = 19991.666
$ws.Range("K136").Value = 55582050
$ws.Range("L136").Value = 59974.99800000001
$ws.Range("M136").Value = -55579500
$ws.Range("N136").Value = -65074.99800000001

